$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header strings: "<name>_old" -> "<name>_FV2404" and "<name>_new" -> "<name>_FV2410"
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2404"
}

# Column K (11) stays "diff"

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2410"
}

# 2. Freeze the top header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# 3. Convert the range into an Excel Table ("Table1")
$tableRange = $ws.Range("A1:U59")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Restore the active selection to A1 (matches original workbook default)
$ws.Range("A1").Select() | Out-Null
